$wb = $excel.ActiveWorkbook

# Update worksheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9559
$ws1.Range("F3").Value = 211
$ws1.Range("F4").Value = 28
$ws1.Range("F5").Value = 534

# Update worksheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9559
$ws4.Range("F3").Value = 211
$ws4.Range("F4").Value = 28
$ws4.Range("F5").Value = 534
